# "changes in Expenses incorporated"
# Adds four new Reimbursement-Form test-case rows (6-9) to the
# TestCaseMaster sheet, renames the shared "FileName" value used by the
# Reimbursement Unit rows (and the new Reimbursement Form rows) from
# "ReimbUnits//ReimbUnitsCreationScenarios.xlsx" to "TestData.xlsx",
# narrows column F, and moves the active selection to E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the shared FileName value for the existing Reimbursement Unit
#     rows first, so every row that used the old file name now shares the
#     new "TestData.xlsx" text (rows 2-5).
$ws.Range("F2").Value = "TestData.xlsx"
$ws.Range("F3").Value = "TestData.xlsx"
$ws.Range("F4").Value = "TestData.xlsx"
$ws.Range("F5").Value = "TestData.xlsx"

# --- Row 6: Reimbursement Form Creation
$ws.Range("B6").Value = "Reimbursement"
$ws.Range("C6").Value = "Reimbursement"
$ws.Range("D6").Value = "Reimbursement Form Creation"
$ws.Range("E6").Value = "com.darwinbox.reimbursement.TestCreateReimbForm"
$ws.Range("F6").Value = "TestData.xlsx"
$ws.Range("G6").Value = "ReimbForm"
$ws.Range("H6").Value = "all"

# --- Row 7: Reimbursement Form Updation
$ws.Range("B7").Value = "Reimbursement"
$ws.Range("C7").Value = "Reimbursement"
$ws.Range("D7").Value = "Reimbursement Form Updation"
$ws.Range("E7").Value = "com.darwinbox.reimbursement.TestUpdateReimbForm"
$ws.Range("F7").Value = "TestData.xlsx"
$ws.Range("G7").Value = "ReimbForm"
$ws.Range("H7").Value = "all"

# --- Row 8: Reimbursement Form Deletion
$ws.Range("B8").Value = "Reimbursement"
$ws.Range("C8").Value = "Reimbursement"
$ws.Range("D8").Value = "Reimbursement Form Deletion"
$ws.Range("E8").Value = "com.darwinbox.reimbursement.TestDeleteReimbForm"
$ws.Range("F8").Value = "TestData.xlsx"
$ws.Range("G8").Value = "ReimbForm"
$ws.Range("H8").Value = "all"

# --- Row 9: Reimbursement Form Duplicate check
$ws.Range("B9").Value = "Reimbursement"
$ws.Range("C9").Value = "Reimbursement"
$ws.Range("D9").Value = "Reimbursement Form Duplicate check"
$ws.Range("E9").Value = "com.darwinbox.reimbursement.TestValidateDuplicateReimbForm"
$ws.Range("F9").Value = "TestData.xlsx"
$ws.Range("G9").Value = "ReimbForm"
$ws.Range("H9").Value = "all"

# --- Column F is narrower now that it only holds "TestData.xlsx".
$ws.Columns.Item(6).ColumnWidth = 15.5

# --- Move the active selection to E6.
$ws.Range("E6").Select() | Out-Null
